# Update average_county_temperature (col I) with NOAA data, and the
# dependent worst/best ASHP COP columns (N, O) that were recalculated
# from it, for the affected facility rows (NAICS 311230 rows re-added
# as part of the merged dataset).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("I7").Value = 19.79629629629628
$ws.Range("N7").Value = 1.383082880591839
$ws.Range("O7").Value = 1.45851929478486
$ws.Range("I8").Value = 16.86342592592595
$ws.Range("N8").Value = 1.358217807733239
$ws.Range("O8").Value = 1.430585745597461
$ws.Range("I9").Value = 5.486111111111112
$ws.Range("N9").Value = 1.26966971746916
$ws.Range("O9").Value = 1.331650485436893
$ws.Range("I10").Value = 14.96875
$ws.Range("N10").Value = 1.34262450293505
$ws.Range("O10").Value = 1.413102197137674
$ws.Range("I11").Value = 14.96875
$ws.Range("I12").Value = 17.25771604938272
$ws.Range("N12").Value = 1.361508482130158
$ws.Range("O12").Value = 1.434278670802308
$ws.Range("I17").Value = 13.75752314814816
$ws.Range("N17").Value = 1.3328422686908
$ws.Range("O17").Value = 1.402147510806076
$ws.Range("I18").Value = 13.75752314814816
$ws.Range("I20").Value = 14.96875
$ws.Range("N20").Value = 1.34262450293505
$ws.Range("O20").Value = 1.413102197137674
$ws.Range("I21").Value = 14.96875
$ws.Range("I24").Value = 13.0158303464755
$ws.Range("N24").Value = 1.326922189449132
$ws.Range("O24").Value = 1.395522876500952
$ws.Range("I25").Value = 13.0158303464755
$ws.Range("I26").Value = 16.86342592592595
$ws.Range("N26").Value = 1.358217807733239
$ws.Range("O26").Value = 1.430585745597461
$ws.Range("I27").Value = -3.847222222222223
$ws.Range("N27").Value = 1.205212661479187
$ws.Range("O27").Value = 1.260158184868579
$ws.Range("I30").Value = 12.41429539295394
$ws.Range("N30").Value = 1.322159312571127
$ws.Range("O30").Value = 1.390195897186759
$ws.Range("I31").Value = 19.60879629629628
$ws.Range("N31").Value = 1.381466033569542
$ws.Range("O31").Value = 1.45670088811587
